$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number must be forced to stay
# text (matching the original inlineStr cells) by temporarily switching the
# cell to a Text number format, then restoring the default "Normal" style so
# no stray style index is left behind on the cell.

$ws.Range('D2').Value = '35.768.93'
$ws.Range('E2').Value = '  -2.55%  '
$ws.Range('D3').Value = '1.989.17'
$ws.Range('E3').Value = '  -3.48%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.78'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.35%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.633'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.10%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '59.57'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +7.61%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '58.94'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.11%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.367'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.65%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0740'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.81%  '
$ws.Range('E12').Value = '  -2.35%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.947'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.16%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '14.74'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.84%  '
$ws.Range('D15').Value = '2.280.67'
$ws.Range('E15').Value = '  -3.43%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.33'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.98%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '19.31'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +11.31%  '
$ws.Range('D18').Value = '1.965.52'
$ws.Range('E18').Value = '  -4.34%  '
$ws.Range('D19').Value = '35.748.80'
$ws.Range('E19').Value = '  -2.41%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.85'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.87%  '
$ws.Range('D21').Value = '0.0₃0850'
$ws.Range('E21').Value = '  -1.85%  '
$ws.Range('E22').Value = '  -1.43%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '233.52'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.24%  '
$ws.Range('E24').Value = '  +0.17%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.59'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +12.70%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.27'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -5.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.55'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.81%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '165.00'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.70%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.33'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.20%  '
$ws.Range('E30').Value = '  -2.37%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.93'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.34%  '
$ws.Range('E32').Value = '  -6.65%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0977'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +13.83%  '
$ws.Range('E34').Value = '  +0.62%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.47'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +7.88%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.40'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.02%  '
$ws.Range('E37').Value = '  +0.14%  '
$ws.Range('E38').Value = '  -1.98%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.72'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +12.20%  '
$ws.Range('E40').Value = '  -1.78%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0945'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.72%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0214'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.17%  '
$ws.Range('B43').Value = 'HuobiToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.85'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.81%  '
$ws.Range('E44').Value = '  -1.02%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '94.17'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.27%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.82'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.54%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '16.41'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.40%  '
$ws.Range('D48').Value = '1.370.27'
$ws.Range('E48').Value = '  -3.49%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.90'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.33%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.33'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.83%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '46.90'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.93%  '
